# Santos is no longer in the top flight; replace the SANTOS row with
# ATHLETICO PR's stats, and move the active selection to F8 (matching the
# cell the author was looking at when making the correction).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Times-br")

# Row 8 held SANTOS (13 estaduais, 9 nacionais, 4 continentais). Overwrite it
# in place with ATHLETICO PR's figures; column E recalculates automatically
# via the shared SUM formula already in the sheet.
$ws.Range("A8").Value = "ATHLETICO PR"
$ws.Range("B8").Value = 28
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 2

# Move the selection to where the author left it after the edit.
$ws.Range("F8").Select() | Out-Null
